$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (O) to the table, carrying over the same
# formatting used by the existing "2020" column (N) for the data rows.
$ws.Range("N4:N14").Copy($ws.Range("O4:O14"))

# The 2021 figures differ from the 2020 ones for a handful of rows;
# overwrite just those values (same "-" placeholder string elsewhere).
$ws.Range("O4").Value = 2021
$ws.Range("O7").Value = 1
$ws.Range("O10").Value = "-"
$ws.Range("O12").Value = 1
$ws.Range("O13").Value = "-"

# Reflect the new active cell/selection shown in the sheet view.
$ws.Range("P1").Select()
